$d = $word.ActiveDocument

# The first paragraph currently reads "Basic if demonstration :". We need
# to splice in an M2Doc validation warning right after "Basic " and before
# "if", consisting of: 4 spaces, "<---", the warning message, and 4 more
# spaces - with "<---" and the message carrying the validation-error
# look (orange, bigger, light-gray highlight).

$findRange = $d.Range(0, 0)
$found = $findRange.Find.Execute("Basic ", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Basic ' at the start of the document."
}

$insertPoint = $findRange.End

$arrow = "<---"
$message = "M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0"
$fullInsertion = "    " + $arrow + $message + "    "

# Insert all the new text in one go (plain formatting, inherited from the
# "Basic " run that precedes it) so nothing-after stays unformatted.
$insertRange = $d.Range($insertPoint, $insertPoint)
$insertRange.InsertAfter($fullInsertion)

# Now go back and re-style just the "<---" + message portion to look like
# the other M2Doc validation errors in this document: orange text, 16pt,
# light-gray highlight.
$warnStart = $insertPoint + 4
$warnEnd = $warnStart + $arrow.Length + $message.Length
$warnRange = $d.Range($warnStart, $warnEnd)
$warnRange.Font.Color = 42495
$warnRange.Font.Size = 16
$warnRange.Font.HighlightColorIndex = 16
